$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21 (shifts existing rows 21-26 down to 22-27,
# preserving their original A/B/C/D/E contents).
$ws.Rows.Item(21).Insert()

# Fill in the brand-new row 21 (event 580 - Mortalidad por dengue).
# The "evento" column stores codes as text (matches the rest of column A),
# so force text formatting before entering the numeric-looking code, then
# drop the formatting override again so no stray style is left behind.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "580"
$ws.Range("A21").ClearFormats()
$ws.Range("B21").Value = "Mortalidad por dengue"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 1

# Update Esperado (C) / Observado (D) / valor p (E) values that changed
# across the rest of the sheet (semana 05 de 2025 refresh).

# Row 2 - 113 Desnutricion aguda en menores de 5 anos
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 0.02

# Row 3 - 115 Cancer en menores de 18 anos
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 0

# Row 4 - 155 Cancer de la mama y cuello uterino
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 0.14

# Row 5 - 210 Dengue
$ws.Range("D5").Value = 45

# Row 6 - 215 Defectos congenitos
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 3

# Row 7 - 220 Dengue grave
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0

# Row 8 - 300 Agresiones por animales potencialmente transmisores de rabia
$ws.Range("C8").Value = 47
$ws.Range("D8").Value = 49
$ws.Range("E8").Value = 0.05

# Row 11 - 346 Ira por virus nuevo
$ws.Range("C11").Value = 154
$ws.Range("D11").Value = 0

# Row 12 - 348 Infeccion respiratoria aguda grave irag inusitada
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = 1

# Row 14 - 356 Intento de suicidio
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 12
$ws.Range("E14").Value = 0.01

# Row 15 - 357 Iad - infecciones asociadas a dispositivos - individual
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.05

# Row 16 - 365 Intoxicaciones
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 0.16

# Row 17 - 455 Leptospirosis
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0.18

# Row 18 - 465 Malaria
$ws.Range("D18").Value = 2

# Row 19 - 549 Morbilidad materna extrema
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 0.04

# Row 20 - 560 Mortalidad perinatal y neonatal tardia
$ws.Range("D20").Value = 1

# Row 22 - 620 Parotiditis (previously row 21)
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0.14

# Row 23 - 740 Sifilis congenita (previously row 22)
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0.37

# Row 24 - 750 Sifilis gestacional (previously row 23)
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 0.1

# Row 25 - 813 Tuberculosis (previously row 24)
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 0.09

# Row 26 - 831 Varicela individual (previously row 25)
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 0.14

# Row 27 - 850 Vih/sida/mortalidad por sida (previously row 26)
$ws.Range("C27").Value = 9
$ws.Range("D27").Value = 7
$ws.Range("E27").Value = 0.12

Write-Host "Done applying semana 05 de 2025 updates"
